$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.855.84'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.23%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.273.57'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.71%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.92%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.38'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.04%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.532'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.69%  '

# Row 8
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
$ws.Range('E9').Value = '  +3.61%  '

# Row 10
$ws.Range('E10').Value = '  +6.48%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.10'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.22%  '

# Row 12
$ws.Range('E12').Value = '  +2.39%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.116'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.48%  '

# Row 14
$ws.Range('E14').Value = '  +3.51%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.625.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.88%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.73%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.278.69'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.09%  '

# Row 18
$ws.Range('E18').Value = '  +3.32%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.795.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.30%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.38%  '

# Row 21
$ws.Range('E21').Value = '  +2.08%  '

# Row 22
$ws.Range('E22').Value = '  +2.59%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.34'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.14%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '243.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.66%  '

# Row 25
$ws.Range('E25').Value = '  +3.54%  '

# Row 26
$ws.Range('E26').Value = '  +0.07%  '

# Row 27
$ws.Range('E27').Value = '  +4.52%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.33'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.08%  '

# Row 29
$ws.Range('E29').Value = '  +3.12%  '

# Row 30
$ws.Range('E30').Value = '  -3.74%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.20%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.92'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.28%  '

# Row 33
$ws.Range('E33').Value = '  -0.02%  '

# Row 34
$ws.Range('E34').Value = '  +4.12%  '

# Row 35
$ws.Range('E35').Value = '  +4.74%  '

# Row 36
$ws.Range('E36').Value = '  -0.12%  '

# Row 37
$ws.Range('E37').Value = '  +3.73%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.50%  '

# Row 39
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.116'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.31%  '

# Row 40
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.105'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.92%  '

# Row 41
$ws.Range('E41').Value = '  +3.04%  '

# Row 42
$ws.Range('E42').Value = '  +5.56%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.084.69'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.53%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.60'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.35%  '

# Row 45
$ws.Range('E45').Value = '  +3.15%  '

# Row 46
$ws.Range('E46').Value = '  +1.95%  '

# Row 47
$ws.Range('E47').Value = '  +6.46%  '

# Row 48
$ws.Range('E48').Value = '  +4.55%  '

# Row 49
$ws.Range('E49').Value = '  +3.11%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.08'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.40%  '

# Row 51
$ws.Range('E51').Value = '  +3.41%  '

